$wb = $excel.ActiveWorkbook

# --- Rename the "tâches" sheet to "tâches User Stories" ---
$wsTaches = $wb.Worksheets.Item("tâches")
$wsTaches.Name = "tâches User Stories"

# --- Update the burndown schedule header row (row 9) on "Sprint 0": ---
# the working-day time slots moved from hourly (09h00..14h00) to 90-minute
# steps (09h00, 10h30, 12h00, 13h30, 15h00, 16h30)
$ws3 = $wb.Worksheets.Item("Sprint 0")
$ws3.Cells.Item(9, 7).Value  = "10h30"
$ws3.Cells.Item(9, 8).Value  = "12h00"
$ws3.Cells.Item(9, 9).Value  = "13h30"
$ws3.Cells.Item(9, 10).Value = "15h00"
$ws3.Cells.Item(9, 11).Value = "16h30"

# --- Preserve/refresh the last selection on each sheet ---
# "tâches User Stories" is the active sheet; just move the selection there.
$wsTaches.Range("E36").Select()

# "Sprint 0" selection moved too; switch to it briefly to set the
# selection, then return focus to "tâches User Stories" so the active
# tab stays the same as before.
$ws3.Activate()
$ws3.Range("M13").Select()
$wsTaches.Activate()
